$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new weekly-report rows before the current row 573, pushing the
# existing rows 573:582 down to 576:585 (dimension grows from T582 to T585).
$ws.Range("A573:A575").EntireRow.Insert()

# New row 573: Artic Snow / Primera
$ws.Range("A573").Value = 10
$ws.Range("B573").Value = "Vega Modelo de Temuco"
$ws.Range("C573").Value = "La Araucanía"
$ws.Range("D573").Value = 44656
$ws.Range("E573").Value = 9
$ws.Range("F573").Value = "Fruta"
$ws.Range("G573").Value = 100103
$ws.Range("H573").Value = "Frutos de hueso (carozo)"
$ws.Range("I573").Value = 100103006
$ws.Range("J573").Value = "Nectarín"
$ws.Range("K573").Value = "Artic Snow"
$ws.Range("L573").Value = "Primera"
$ws.Range("M573").Value = 120
$ws.Range("N573").Value = 20000
$ws.Range("O573").Value = 20000
$ws.Range("P573").Value = 20000
$ws.Range("Q573").Value = "$/bandeja 18 kilos granel"
$ws.Range("R573").Value = "Región de O'Higgins"
$ws.Range("S573").Value = 1111
$ws.Range("T573").Value = 18

# New row 574: Artic Snow / Primera (bins)
$ws.Range("A574").Value = 10
$ws.Range("B574").Value = "Vega Modelo de Temuco"
$ws.Range("C574").Value = "La Araucanía"
$ws.Range("D574").Value = 44656
$ws.Range("E574").Value = 9
$ws.Range("F574").Value = "Fruta"
$ws.Range("G574").Value = 100103
$ws.Range("H574").Value = "Frutos de hueso (carozo)"
$ws.Range("I574").Value = 100103006
$ws.Range("J574").Value = "Nectarín"
$ws.Range("K574").Value = "Artic Snow"
$ws.Range("L574").Value = "Primera"
$ws.Range("M574").Value = 2
$ws.Range("N574").Value = 460000
$ws.Range("O574").Value = 460000
$ws.Range("P574").Value = 460000
$ws.Range("Q574").Value = "$/bins (420 kilos)"
$ws.Range("R574").Value = "Región de O'Higgins"
$ws.Range("S574").Value = 1095
$ws.Range("T574").Value = 420

# New row 575: August Red / Especial
$ws.Range("A575").Value = 10
$ws.Range("B575").Value = "Vega Modelo de Temuco"
$ws.Range("C575").Value = "La Araucanía"
$ws.Range("D575").Value = 44656
$ws.Range("E575").Value = 9
$ws.Range("F575").Value = "Fruta"
$ws.Range("G575").Value = 100103
$ws.Range("H575").Value = "Frutos de hueso (carozo)"
$ws.Range("I575").Value = 100103006
$ws.Range("J575").Value = "Nectarín"
$ws.Range("K575").Value = "August Red"
$ws.Range("L575").Value = "Especial"
$ws.Range("M575").Value = 50
$ws.Range("N575").Value = 24000
$ws.Range("O575").Value = 24000
$ws.Range("P575").Value = 24000
$ws.Range("Q575").Value = "$/bandeja 18 kilos granel"
$ws.Range("R575").Value = "Región de O'Higgins"
$ws.Range("S575").Value = 1333
$ws.Range("T575").Value = 18
